$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:I1").UnMerge()
$ws.Columns.Item(1).Insert()

for ($c = 1; $c -le 8; $c++) {
    $v = $ws.Cells.Item(1, $c+1).Value()
    $ws.Cells.Item(1, $c).Value = $v
}
$ws.Cells.Item(1,9).Value = ""

Write-Host "row1 final:"
for ($c = 1; $c -le 9; $c++) {
    Write-Host "c$c :[" $ws.Cells.Item(1, $c).Value() "]"
}

$ws.Range("C1:D1").Merge()
$ws.Range("E1:F1").Merge()
$ws.Range("G1:H1").Merge()
Write-Host "merge fixed"
